# chnage in excel data
# Update the AccountCreationData test-data sheet: refresh the sample
# account rows (new emails / usernames / password) and move the active
# selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

# Row 2
$ws.Range("A2").Value = "newt@gmail.com"
$ws.Range("B2").Value = "TestUsee"
$ws.Range("C2").Value = "UserTest"
$ws.Range("D2").Value = "Test1234"

# Row 3
$ws.Range("A3").Value = "qas2@gmail.com"
$ws.Range("B3").Value = "TestUsew"
$ws.Range("C3").Value = "UserTest"
$ws.Range("D3").Value = "Test1234"

# Row 4
$ws.Range("A4").Value = "qad3@gmail.com"
$ws.Range("B4").Value = "TestUsey"
$ws.Range("C4").Value = "UserTest"
$ws.Range("D4").Value = "Test1234"

# Move selection to B4 on the AccountCreationData sheet (stays the active tab)
$ws.Activate()
$ws.Range("B4").Select()
